# Auto-generated edit script: updates cached market/profit values
# across multiple leve-profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 6835.7896  # ALC!H17: 6120.4546 -> 6835.7896
$ws.Cells.Item(17, 10).Value = 7132.222  # ALC!J17: 6340.476 -> 7132.222
$ws.Cells.Item(17, 12).Value = 21396.666  # ALC!L17: 19021.428 -> 21396.666
$ws.Cells.Item(17, 14).Value = -21732.666  # ALC!N17: -19357.428 -> -21732.666

$ws.Cells.Item(101, 8).Value = 1132.8334  # ALC!H101: 1559.0769 -> 1132.8334
$ws.Cells.Item(101, 9).Value = 917.7778  # ALC!I101: 992.6667 -> 917.7778
$ws.Cells.Item(101, 10).Value = 1778  # ALC!J101: 2833.5 -> 1778
$ws.Cells.Item(101, 11).Value = 2753.3334  # ALC!K101: 2978.0001 -> 2753.3334
$ws.Cells.Item(101, 12).Value = 5334  # ALC!L101: 8500.5 -> 5334
$ws.Cells.Item(101, 13).Value = -1131.3334  # ALC!M101: -1356.0001 -> -1131.3334
$ws.Cells.Item(101, 14).Value = -8578  # ALC!N101: -11744.5 -> -8578

$ws.Cells.Item(116, 8).Value = 5869.7334  # ALC!H116: 6790.7856 -> 5869.7334
$ws.Cells.Item(116, 9).Value = 4279.5557  # ALC!I116: 5692.625 -> 4279.5557
$ws.Cells.Item(116, 11).Value = 4279.5557  # ALC!K116: 5692.625 -> 4279.5557
$ws.Cells.Item(116, 13).Value = -837.5556999999999  # ALC!M116: -2250.625 -> -837.5556999999999

$ws.Cells.Item(137, 8).Value = 3018.484  # ALC!H137: 3053.2131 -> 3018.484
$ws.Cells.Item(137, 9).Value = 2241.0908  # ALC!I137: 2375.2 -> 2241.0908
$ws.Cells.Item(137, 11).Value = 6723.2724  # ALC!K137: 7125.599999999999 -> 6723.2724
$ws.Cells.Item(137, 13).Value = -4173.2724  # ALC!M137: -4575.599999999999 -> -4173.2724

$ws.Cells.Item(138, 8).Value = 2951.923  # ALC!H138: 2937.782 -> 2951.923
$ws.Cells.Item(138, 9).Value = 1794.3  # ALC!I138: 1765.2858 -> 1794.3
$ws.Cells.Item(138, 10).Value = 3351.1035  # ALC!J138: 3369.7544 -> 3351.1035
$ws.Cells.Item(138, 11).Value = 5382.9  # ALC!K138: 5295.857400000001 -> 5382.9
$ws.Cells.Item(138, 12).Value = 10053.3105  # ALC!L138: 10109.2632 -> 10053.3105
$ws.Cells.Item(138, 13).Value = -242.8999999999996  # ALC!M138: -155.8574000000008 -> -242.8999999999996
$ws.Cells.Item(138, 14).Value = -20333.3105  # ALC!N138: -20389.2632 -> -20333.3105

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 5179.346  # ARM!H2: 5213.9614 -> 5179.346
$ws.Cells.Item(2, 9).Value = 833.1739  # ARM!I2: 898.3182 -> 833.1739
$ws.Cells.Item(2, 10).Value = 38500  # ARM!J2: 28950 -> 38500
$ws.Cells.Item(2, 11).Value = 833.1739  # ARM!K2: 898.3182 -> 833.1739
$ws.Cells.Item(2, 12).Value = 38500  # ARM!L2: 28950 -> 38500
$ws.Cells.Item(2, 13).Value = -720.1739  # ARM!M2: -785.3182 -> -720.1739
$ws.Cells.Item(2, 14).Value = -38726  # ARM!N2: -29176 -> -38726

$ws.Cells.Item(74, 8).Value = 11496630  # ARM!H74: 10755024 -> 11496630
$ws.Cells.Item(74, 9).Value = 18520126  # ARM!I74: 17545558 -> 18520126
$ws.Cells.Item(74, 10).Value = 3634.7273  # ARM!J74: 3345.5 -> 3634.7273
$ws.Cells.Item(74, 11).Value = 18520126  # ARM!K74: 17545558 -> 18520126
$ws.Cells.Item(74, 12).Value = 3634.7273  # ARM!L74: 3345.5 -> 3634.7273
$ws.Cells.Item(74, 13).Value = -18519252  # ARM!M74: -17544684 -> -18519252
$ws.Cells.Item(74, 14).Value = -5382.7273  # ARM!N74: -5093.5 -> -5382.7273

$ws.Cells.Item(77, 8).Value = 11496630  # ARM!H77: 10755024 -> 11496630
$ws.Cells.Item(77, 9).Value = 18520126  # ARM!I77: 17545558 -> 18520126
$ws.Cells.Item(77, 10).Value = 3634.7273  # ARM!J77: 3345.5 -> 3634.7273
$ws.Cells.Item(77, 11).Value = 92600630  # ARM!K77: 87727790 -> 92600630
$ws.Cells.Item(77, 12).Value = 18173.6365  # ARM!L77: 16727.5 -> 18173.6365
$ws.Cells.Item(77, 13).Value = -92596262  # ARM!M77: -87723422 -> -92596262
$ws.Cells.Item(77, 14).Value = -26909.6365  # ARM!N77: -25463.5 -> -26909.6365

$ws.Cells.Item(97, 8).Value = 4250  # ARM!H97: 3500 -> 4250
$ws.Cells.Item(97, 10).Value = 5000  # ARM!J97: 0 -> 5000
$ws.Cells.Item(97, 12).Value = 5000  # ARM!L97: 0 -> 5000
$ws.Cells.Item(97, 14).Value = -5992  # ARM!N97: None -> -5992

$ws.Cells.Item(110, 8).Value = 5437  # ARM!H110: 5421.657 -> 5437
$ws.Cells.Item(110, 10).Value = 6401  # ARM!J110: 6275.9165 -> 6401
$ws.Cells.Item(110, 12).Value = 6401  # ARM!L110: 6275.9165 -> 6401
$ws.Cells.Item(110, 14).Value = -10491  # ARM!N110: -10365.9165 -> -10491

$ws.Cells.Item(116, 8).Value = 5179.346  # ARM!H116: 5213.9614 -> 5179.346
$ws.Cells.Item(116, 9).Value = 833.1739  # ARM!I116: 898.3182 -> 833.1739
$ws.Cells.Item(116, 10).Value = 38500  # ARM!J116: 28950 -> 38500
$ws.Cells.Item(116, 11).Value = 833.1739  # ARM!K116: 898.3182 -> 833.1739
$ws.Cells.Item(116, 12).Value = 38500  # ARM!L116: 28950 -> 38500
$ws.Cells.Item(116, 13).Value = 1460.8261  # ARM!M116: 1395.6818 -> 1460.8261
$ws.Cells.Item(116, 14).Value = -43088  # ARM!N116: -33538 -> -43088

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 5179.346  # BSM!H3: 5213.9614 -> 5179.346
$ws.Cells.Item(3, 9).Value = 833.1739  # BSM!I3: 898.3182 -> 833.1739
$ws.Cells.Item(3, 10).Value = 38500  # BSM!J3: 28950 -> 38500
$ws.Cells.Item(3, 11).Value = 833.1739  # BSM!K3: 898.3182 -> 833.1739
$ws.Cells.Item(3, 12).Value = 38500  # BSM!L3: 28950 -> 38500
$ws.Cells.Item(3, 13).Value = -719.1739  # BSM!M3: -784.3182 -> -719.1739
$ws.Cells.Item(3, 14).Value = -38728  # BSM!N3: -29178 -> -38728

$ws.Cells.Item(94, 8).Value = 1624.75  # BSM!H94: 2500 -> 1624.75
$ws.Cells.Item(94, 9).Value = 1624.75  # BSM!I94: 2500 -> 1624.75
$ws.Cells.Item(94, 11).Value = 1624.75  # BSM!K94: 2500 -> 1624.75
$ws.Cells.Item(94, 13).Value = -1173.75  # BSM!M94: -2049 -> -1173.75

$ws.Cells.Item(99, 8).Value = 2301.2  # BSM!H99: 5002.5 -> 2301.2
$ws.Cells.Item(99, 9).Value = 2248.75  # BSM!I99: 5002.5 -> 2248.75
$ws.Cells.Item(99, 10).Value = 2511  # BSM!J99: 0 -> 2511
$ws.Cells.Item(99, 11).Value = 2248.75  # BSM!K99: 5002.5 -> 2248.75
$ws.Cells.Item(99, 12).Value = 2511  # BSM!L99: 0 -> 2511
$ws.Cells.Item(99, 13).Value = -750.75  # BSM!M99: -3504.5 -> -750.75
$ws.Cells.Item(99, 14).Value = -5507  # BSM!N99: None -> -5507

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 26000.682  # CRP!H31: 31317.805 -> 26000.682
$ws.Cells.Item(31, 9).Value = 1661.871  # CRP!I31: 1691.32 -> 1661.871
$ws.Cells.Item(31, 10).Value = 84039.38  # CRP!J31: 98650.73 -> 84039.38
$ws.Cells.Item(31, 11).Value = 1661.871  # CRP!K31: 1691.32 -> 1661.871
$ws.Cells.Item(31, 12).Value = 84039.38  # CRP!L31: 98650.73 -> 84039.38
$ws.Cells.Item(31, 13).Value = -1366.871  # CRP!M31: -1396.32 -> -1366.871
$ws.Cells.Item(31, 14).Value = -84629.38  # CRP!N31: -99240.73 -> -84629.38

$ws.Cells.Item(34, 8).Value = 26000.682  # CRP!H34: 31317.805 -> 26000.682
$ws.Cells.Item(34, 9).Value = 1661.871  # CRP!I34: 1691.32 -> 1661.871
$ws.Cells.Item(34, 10).Value = 84039.38  # CRP!J34: 98650.73 -> 84039.38
$ws.Cells.Item(34, 11).Value = 1661.871  # CRP!K34: 1691.32 -> 1661.871
$ws.Cells.Item(34, 12).Value = 84039.38  # CRP!L34: 98650.73 -> 84039.38
$ws.Cells.Item(34, 13).Value = -1459.871  # CRP!M34: -1489.32 -> -1459.871
$ws.Cells.Item(34, 14).Value = -84443.38  # CRP!N34: -99054.73 -> -84443.38

$ws.Cells.Item(58, 8).Value = 3130.4285  # CRP!H58: 3183.3809 -> 3130.4285
$ws.Cells.Item(58, 9).Value = 1631.0588  # CRP!I58: 1696.4706 -> 1631.0588
$ws.Cells.Item(58, 11).Value = 1631.0588  # CRP!K58: 1696.4706 -> 1631.0588
$ws.Cells.Item(58, 13).Value = -1428.0588  # CRP!M58: -1493.4706 -> -1428.0588

$ws.Cells.Item(62, 8).Value = 7516.5  # CRP!H62: 9488.888999999999 -> 7516.5
$ws.Cells.Item(62, 9).Value = 2998.8572  # CRP!I62: 3723.5 -> 2998.8572
$ws.Cells.Item(62, 10).Value = 13841.2  # CRP!J62: 14101.2 -> 13841.2
$ws.Cells.Item(62, 11).Value = 2998.8572  # CRP!K62: 3723.5 -> 2998.8572
$ws.Cells.Item(62, 12).Value = 13841.2  # CRP!L62: 14101.2 -> 13841.2
$ws.Cells.Item(62, 13).Value = -2374.8572  # CRP!M62: -3099.5 -> -2374.8572
$ws.Cells.Item(62, 14).Value = -15089.2  # CRP!N62: -15349.2 -> -15089.2

$ws.Cells.Item(65, 8).Value = 7516.5  # CRP!H65: 9488.888999999999 -> 7516.5
$ws.Cells.Item(65, 9).Value = 2998.8572  # CRP!I65: 3723.5 -> 2998.8572
$ws.Cells.Item(65, 10).Value = 13841.2  # CRP!J65: 14101.2 -> 13841.2
$ws.Cells.Item(65, 11).Value = 14994.286  # CRP!K65: 18617.5 -> 14994.286
$ws.Cells.Item(65, 12).Value = 69206  # CRP!L65: 70506 -> 69206
$ws.Cells.Item(65, 13).Value = -11874.286  # CRP!M65: -15497.5 -> -11874.286
$ws.Cells.Item(65, 14).Value = -75446  # CRP!N65: -76746 -> -75446

$ws.Cells.Item(125, 8).Value = 0  # CRP!H125: 49000 -> 0
$ws.Cells.Item(125, 10).Value = 0  # CRP!J125: 49000 -> 0
$ws.Cells.Item(125, 12).Value = 0  # CRP!L125: 49000 -> 0
$ws.Cells.Item(125, 14).ClearContents()  # CRP!N125: remove (was -53920)

$ws.Cells.Item(132, 8).Value = 3261.2144  # CRP!H132: 3542.261 -> 3261.2144
$ws.Cells.Item(132, 9).Value = 2652.7058  # CRP!I132: 2937.8333 -> 2652.7058
$ws.Cells.Item(132, 11).Value = 7958.117400000001  # CRP!K132: 8813.499899999999 -> 7958.117400000001
$ws.Cells.Item(132, 13).Value = -5428.117400000001  # CRP!M132: -6283.499899999999 -> -5428.117400000001

$ws.Cells.Item(134, 8).Value = 3143.7917  # CRP!H134: 3192.9048 -> 3143.7917
$ws.Cells.Item(134, 9).Value = 2054.353  # CRP!I134: 1894.5714 -> 2054.353
$ws.Cells.Item(134, 11).Value = 6163.059  # CRP!K134: 5683.7142 -> 6163.059
$ws.Cells.Item(134, 13).Value = -3628.059  # CRP!M134: -3148.7142 -> -3628.059

$ws.Cells.Item(136, 8).Value = 3130.4285  # CRP!H136: 3183.3809 -> 3130.4285
$ws.Cells.Item(136, 9).Value = 1631.0588  # CRP!I136: 1696.4706 -> 1631.0588
$ws.Cells.Item(136, 11).Value = 4893.1764  # CRP!K136: 5089.4118 -> 4893.1764
$ws.Cells.Item(136, 13).Value = -2343.1764  # CRP!M136: -2539.4118 -> -2343.1764

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 63.8  # CUL!H12: 80 -> 63.8
$ws.Cells.Item(12, 9).Value = 10  # CUL!I12: 0 -> 10
$ws.Cells.Item(12, 10).Value = 77.25  # CUL!J12: 80 -> 77.25
$ws.Cells.Item(12, 11).Value = 30  # CUL!K12: 0 -> 30
$ws.Cells.Item(12, 12).Value = 231.75  # CUL!L12: 240 -> 231.75
$ws.Cells.Item(12, 13).Value = 143  # CUL!M12: None -> 143
$ws.Cells.Item(12, 14).Value = -577.75  # CUL!N12: -586 -> -577.75

$ws.Cells.Item(81, 8).Value = 4183.6665  # CUL!H81: 2827.6365 -> 4183.6665
$ws.Cells.Item(81, 9).Value = 1273.625  # CUL!I81: 1398.7778 -> 1273.625
$ws.Cells.Item(81, 10).Value = 10003.75  # CUL!J81: 9257.5 -> 10003.75
$ws.Cells.Item(81, 11).Value = 3820.875  # CUL!K81: 4196.3334 -> 3820.875
$ws.Cells.Item(81, 12).Value = 30011.25  # CUL!L81: 27772.5 -> 30011.25
$ws.Cells.Item(81, 13).Value = -2697.875  # CUL!M81: -3073.3334 -> -2697.875
$ws.Cells.Item(81, 14).Value = -32257.25  # CUL!N81: -30018.5 -> -32257.25

$ws.Cells.Item(84, 8).Value = 4183.6665  # CUL!H84: 2827.6365 -> 4183.6665
$ws.Cells.Item(84, 9).Value = 1273.625  # CUL!I84: 1398.7778 -> 1273.625
$ws.Cells.Item(84, 10).Value = 10003.75  # CUL!J84: 9257.5 -> 10003.75
$ws.Cells.Item(84, 11).Value = 11462.625  # CUL!K84: 12589.0002 -> 11462.625
$ws.Cells.Item(84, 12).Value = 90033.75  # CUL!L84: 83317.5 -> 90033.75
$ws.Cells.Item(84, 13).Value = -5846.625  # CUL!M84: -6973.0002 -> -5846.625
$ws.Cells.Item(84, 14).Value = -101265.75  # CUL!N84: -94549.5 -> -101265.75

$ws.Cells.Item(131, 8).Value = 11410466  # CUL!H131: 11410445 -> 11410466
$ws.Cells.Item(131, 10).Value = 8455973  # CUL!J131: 8455947 -> 8455973
$ws.Cells.Item(131, 12).Value = 25367919  # CUL!L131: 25367841 -> 25367919
$ws.Cells.Item(131, 14).Value = -25377999  # CUL!N131: -25377921 -> -25377999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(47, 8).Value = 1200000  # GSM!H47: 1100000 -> 1200000
$ws.Cells.Item(47, 10).Value = 1200000  # GSM!J47: 1100000 -> 1200000
$ws.Cells.Item(47, 12).Value = 1200000  # GSM!L47: 1100000 -> 1200000
$ws.Cells.Item(47, 14).Value = -1201136  # GSM!N47: -1101136 -> -1201136

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2686.1155  # LTW!H22: 2450.3667 -> 2686.1155
$ws.Cells.Item(22, 9).Value = 1211.6875  # LTW!I22: 1180.1177 -> 1211.6875
$ws.Cells.Item(22, 10).Value = 5045.2  # LTW!J22: 4111.4614 -> 5045.2
$ws.Cells.Item(22, 11).Value = 1211.6875  # LTW!K22: 1180.1177 -> 1211.6875
$ws.Cells.Item(22, 12).Value = 5045.2  # LTW!L22: 4111.4614 -> 5045.2
$ws.Cells.Item(22, 13).Value = -916.6875  # LTW!M22: -885.1177 -> -916.6875
$ws.Cells.Item(22, 14).Value = -5635.2  # LTW!N22: -4701.4614 -> -5635.2

$ws.Cells.Item(27, 8).Value = 2686.1155  # LTW!H27: 2450.3667 -> 2686.1155
$ws.Cells.Item(27, 9).Value = 1211.6875  # LTW!I27: 1180.1177 -> 1211.6875
$ws.Cells.Item(27, 10).Value = 5045.2  # LTW!J27: 4111.4614 -> 5045.2
$ws.Cells.Item(27, 11).Value = 1211.6875  # LTW!K27: 1180.1177 -> 1211.6875
$ws.Cells.Item(27, 12).Value = 5045.2  # LTW!L27: 4111.4614 -> 5045.2
$ws.Cells.Item(27, 13).Value = -1104.6875  # LTW!M27: -1073.1177 -> -1104.6875
$ws.Cells.Item(27, 14).Value = -5259.2  # LTW!N27: -4325.4614 -> -5259.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(23, 8).Value = 1799.6666  # WVR!H23: 1966.3334 -> 1799.6666
$ws.Cells.Item(23, 10).Value = 1499  # WVR!J23: 1999 -> 1499
$ws.Cells.Item(23, 12).Value = 1499  # WVR!L23: 1999 -> 1499
$ws.Cells.Item(23, 14).Value = -1957  # WVR!N23: -2457 -> -1957

$ws.Cells.Item(58, 8).Value = 18968.9  # WVR!H58: 18243.223 -> 18968.9

$ws.Cells.Item(82, 8).Value = 27062  # WVR!H82: 29749.666 -> 27062
$ws.Cells.Item(82, 10).Value = 27062  # WVR!J82: 29749.666 -> 27062
$ws.Cells.Item(82, 12).Value = 27062  # WVR!L82: 29749.666 -> 27062
$ws.Cells.Item(82, 14).Value = -27828  # WVR!N82: -30515.666 -> -27828

$ws.Cells.Item(85, 8).Value = 27062  # WVR!H85: 29749.666 -> 27062
$ws.Cells.Item(85, 10).Value = 27062  # WVR!J85: 29749.666 -> 27062
$ws.Cells.Item(85, 12).Value = 27062  # WVR!L85: 29749.666 -> 27062
$ws.Cells.Item(85, 14).Value = -29714  # WVR!N85: -32401.666 -> -29714

$ws.Cells.Item(131, 8).Value = 88209.336  # WVR!H131: 0 -> 88209.336
$ws.Cells.Item(131, 10).Value = 88209.336  # WVR!J131: 0 -> 88209.336
$ws.Cells.Item(131, 12).Value = 88209.336  # WVR!L131: 0 -> 88209.336
$ws.Cells.Item(131, 14).Value = -98289.336  # WVR!N131: None -> -98289.336

$ws.Cells.Item(132, 8).Value = 1509.5769  # WVR!H132: 1509 -> 1509.5769
$ws.Cells.Item(132, 9).Value = 1146.2046  # WVR!I132: 1139.5111 -> 1146.2046
$ws.Cells.Item(132, 10).Value = 3508.125  # WVR!J132: 3884.2856 -> 3508.125
$ws.Cells.Item(132, 11).Value = 3438.6138  # WVR!K132: 3418.5333 -> 3438.6138
$ws.Cells.Item(132, 12).Value = 10524.375  # WVR!L132: 11652.8568 -> 10524.375
$ws.Cells.Item(132, 13).Value = -908.6138000000001  # WVR!M132: -888.5333000000001 -> -908.6138000000001
$ws.Cells.Item(132, 14).Value = -15584.375  # WVR!N132: -16712.8568 -> -15584.375

Write-Output "edit.ps1 applied: 194 cell updates across 8 sheets"
